$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain text that often looks numeric (it uses
# "." as a thousands separator in some rows, e.g. "30.521.96", and keeps
# trailing zeros, e.g. "5.220"). Force each Price cell we touch to Text format
# *before* writing its new value so Excel stores the literal string instead of
# silently converting it to a floating point number.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.521.96"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "1.918.10"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "238.87"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "0.4793"
$ws.Range("E7").Value = "  -1.72%  "
$ws.Range("D8").Value = "0.2869"
$ws.Range("E8").Value = "  -3.34%  "
$ws.Range("D9").Value = "0.06681"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "18.73"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").Value = "103.69"
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.933.67"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07742"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "5.220"
$ws.Range("E14").Value = "  -5.07%  "
$ws.Range("D15").Value = "0.6793"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").Value = "265.52"
$ws.Range("E16").Value = "  -6.63%  "
$ws.Range("D17").Value = "30.565.89"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "0.000007530"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").Value = "12.69"
$ws.Range("E20").Value = "  -4.30%  "
$ws.Range("D21").Value = "5.420"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "6.317"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("D24").Value = "9.611"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "163.09"
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("D26").Value = "18.99"
$ws.Range("E26").Value = "  -5.16%  "
$ws.Range("D27").Value = "2.093"
$ws.Range("E27").Value = "  -5.82%  "
$ws.Range("D28").Value = "0.1022"
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("D29").Value = "1.387"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "4.526"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("D31").Value = "1.513"
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("D32").Value = "4.244"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("D33").Value = "0.04749"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("D34").Value = "0.7356"
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("D35").Value = "1.120"
$ws.Range("E35").Value = "  -4.68%  "
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "2.683"
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("D38").Value = "0.01936"
$ws.Range("E38").Value = "  -4.45%  "
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "6.325"
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("D41").Value = "74.98"
$ws.Range("E41").Value = "  -4.31%  "
$ws.Range("D42").Value = "2.004"
$ws.Range("E42").Value = "  -5.81%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "106.36"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8589"
$ws.Range("E44").Value = "  -5.37%  "
$ws.Range("D45").Value = "0.4280"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").Value = "1.000.59"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "7.495"
$ws.Range("E48").Value = "  -8.47%  "
$ws.Range("D49").Value = "0.1203"
$ws.Range("E49").Value = "  -4.84%  "
$ws.Range("D50").Value = "35.13"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").Value = "8.972"
$ws.Range("E51").Value = "  -4.27%  "
